$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 53, shifting existing rows 53:74 down to 54:75.
$ws.Rows("53:53").Insert()

# Populate the newly inserted row 53 with the new weekly entry.
$ws.Cells.Item(53, 1).Value = 7
$ws.Cells.Item(53, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(53, 3).Value = "Ñuble"
$ws.Cells.Item(53, 4).Value = 44609
$ws.Cells.Item(53, 5).Value = 16
$ws.Cells.Item(53, 6).Value = "Fruta"
$ws.Cells.Item(53, 7).Value = 100108
$ws.Cells.Item(53, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(53, 9).Value = 100108002
$ws.Cells.Item(53, 10).Value = "Mango"
$ws.Cells.Item(53, 11).Value = "Sin especificar"
$ws.Cells.Item(53, 12).Value = "Primera"
$ws.Cells.Item(53, 13).Value = 120
$ws.Cells.Item(53, 14).Value = 7000
$ws.Cells.Item(53, 15).Value = 7500
$ws.Cells.Item(53, 16).Value = 7250
$ws.Cells.Item(53, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(53, 18).Value = "Perú"
$ws.Cells.Item(53, 19).Value = 1812
$ws.Cells.Item(53, 20).Value = 4
